# Journal de travail - reformulation des questions sur la gestion des entretiens
$wb = $excel.ActiveWorkbook
$journal = $wb.Worksheets.Item("Journal")
$totaux  = $wb.Worksheets.Item("Totaux")

# --- 1) Append the 4 new Journal rows (35-38) --------------------------------
# Copy the formatting of the last existing data row (A34:E34) onto the new
# rows so date/number formats & alignment match the rest of the table.
$journal.Range("A34:E34").Copy()
$journal.Range("A35:E38").PasteSpecial(-4122)

$journal.Range("A35").Value2 = 44987
$journal.Range("B35").Value2 = 4
$journal.Range("C35").Value2 = 0.041666666666666664
$journal.Range("D35").Value2 = "Meeting"
$journal.Range("E35").Value2 = "Séance de review avec les personnes d'Eldora"

# --- 2) Rename the "Commentaire" column to "Commentaire/Remarque" -----------
# (updates both the header cell F1 and the Tableau1 column name)
$journal.Range("F1").Value2 = "Commentaire/Remarque"

$journal.Range("A36").Value2 = 44987
$journal.Range("B36").Value2 = 4
$journal.Range("C36").Value2 = 0.020833333333333332
$journal.Range("D36").Value2 = "Documentation"
$journal.Range("E36").Value2 = "Rédaction du rapport de projet"

$journal.Range("A37").Value2 = 44988
$journal.Range("B37").Value2 = 4
$journal.Range("C37").Value2 = 0.07291666666666667
$journal.Range("D37").Value2 = "Documentation"
$journal.Range("E37").Value2 = "Reformulation de toutes les questions sur la gestion des entretiens"

$journal.Range("A38").Value2 = 44988
$journal.Range("B38").Value2 = 4
$journal.Range("C38").Value2 = 0.020833333333333332
$journal.Range("D38").Value2 = "Documentation"
$journal.Range("E38").Value2 = "Rédaction du rapport de projet"

# --- 3) Append the 2 new weekly-summary rows on "Totaux" ---------------------
# Push the previous totals row (old row 13) down to row 15, then fill the new
# weekly rows 13 & 14, copying the formatting from the old totals row first.
$totaux.Range("A13:B13").Copy()
$totaux.Range("A13:B15").PasteSpecial(-4122)

$totaux.Range("A13").Value2 = 44987
$totaux.Range("B13").Formula = "=SUM(Journal!C35:C36)"

$totaux.Range("A14").Value2 = 44988
$totaux.Range("B14").Formula = "=SUM(Journal!C37:C38)"

$totaux.Range("A15").Value2 = "Total"
$totaux.Range("B15").Formula = "=SUM(B2:B14)"

# --- 4) Fix up view/selection state to match the edited workbook ------------
$journal.Range("F38").Select()
$journal.Application.ActiveWindow.ScrollRow = 7

$totaux.Range("D20").Select()
